$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Previous closing" (B), "Closing" (C) and "Change (%)" (D) values
# for each sector row (rows 2-12). Column E (Year to Date Change) is unchanged.

$updates = @(
    @{ Row = 2;  B = 107.33; C = 107.19; D = -0.13 },
    @{ Row = 3;  B = 174.28; C = 173.58; D = -0.4  },
    @{ Row = 4;  B = 993.28; C = 943.92; D = -4.97 },
    @{ Row = 5;  B = 213.8;  C = 213.43; D = -0.17 },
    @{ Row = 6;  B = 313.06; C = 313.91; D = 0.27  },
    @{ Row = 7;  B = 86.82;  C = 86.37;  D = -0.52 },
    @{ Row = 8;  B = 102.91; C = 101.83; D = -1.05 },
    @{ Row = 9;  B = 99.72;  C = 99.55;  D = -0.17 },
    @{ Row = 10; B = 104.02; C = 103.77; D = -0.24 },
    @{ Row = 11; B = 524.94; C = 526.04; D = 0.21  },
    @{ Row = 12; B = 373.92; C = 377.5;  D = 0.96  }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
